$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update G2 and I2 to use a formula (60*48 = 2880) instead of the static 1440 value
$ws.Range("G2").Formula = "=60*48"
$ws.Range("I2").Formula = "=60*48"

# Update the frozen pane / view: split moves from column B to column D,
# and the selection moves from B8:B10 (active B8) to a single cell I16
$ws.Range("D1").Select()
$excel.ActiveWindow.FreezePanes = $true

$ws.Range("I16").Select()
